$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the "EXE Regisration " sheet, placing the copy right after it,
# then rename the new sheet to "EXE Login".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "EXE Login"

# --- Update the Login sheet content ---
# Module column (G) for every test row: "Registration Module" -> "Login Module"
$ws2.Range("G5:G12").Value = "Login Module"

# Row 5 becomes the "Login" test case instead of the "Registration" one.
$ws2.Range("E5").Value = "Verify user can login with valid email and password"
$ws2.Range("I5").Value = "Login with valid email and password"
$ws2.Range("J5").Value = "User successfully logged in"
$ws2.Range("K5").Value = "User successfully logged in"

# Test date for row 5 moves forward a day, matching the rest of the sheet.
$ws2.Range("P5").Value = 46077

# Column J (10) widens slightly to fit the new text (best-fit column).
$ws2.Columns.Item(10).ColumnWidth = 23.6

# --- View/selection state ---
# Original sheet: no longer the selected tab; scrolled down, D13 selected.
$ws1.Activate()
$ws1.Range("D13").Select() | Out-Null

# New Login sheet: becomes the selected tab; K5 selected.
$ws2.Activate()
$ws2.Range("K5").Select() | Out-Null
